$d = $word.ActiveDocument

# Locate the paragraph that starts the "Exibição de Localização GPS" bullet
# item and remove it in its entirety (including its paragraph mark), so the
# list collapses from "Traçar Rotas Principais" straight to "Responsividade
# para Smartphones".
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Exibição de Localização GPS*") {
        $p.Range.Delete()
        break
    }
}
